# Update cryptos list (prices / 1h volume %, plus a few reordered coin rows)
# Prices in column D that look like plain numbers are written with a leading
# apostrophe (then style reset to "Normal") so Excel keeps them as text and
# doesn't strip significant trailing zeros (e.g. "1.00", "8.10").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.826.78"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "3.127.50"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'598.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.27%  "
$ws.Range("D6").Value = "'139.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.59%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.126.39"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("E10").Value = "  -3.17%  "
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "'34.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.60%  "
$ws.Range("D15").Value = "3.641.81"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("E16").Value = "  +2.42%  "
$ws.Range("D17").Value = "63.747.85"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "3.131.37"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "'6.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'480.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.95%  "
$ws.Range("D22").Value = "'0.703"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").Value = "'7.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("D24").Value = "'87.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.20%  "
$ws.Range("D25").Value = "'13.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.01%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'8.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.16%  "
$ws.Range("D29").Value = "'6.95"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.40%  "
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").Value = "'27.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.63%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("D34").Value = "'2.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.65%  "
$ws.Range("D35").Value = "'1.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.43%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'52.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("D38").Value = "0.0₃0728"
$ws.Range("E38").Value = "  -7.98%  "
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.29%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'427.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.23%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("D43").Value = "'8.27"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "2.884.69"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").Value = "'0.257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.70%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -3.81%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "'25.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.19%  "
$ws.Range("D51").Value = "'120.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.52%  "
